{"js": "// Update the date line and every \"A\u00d7B=C\" answer cell to the new values\n// from the commit. Each old string is unique in the document, so a\n// matchCase exact search + in-place replace is sufficient and keeps the\n// existing run formatting (fonts/size) untouched.\nconst replacements = [\n  [\"2025-03-21 Friday\", \"2025-03-22 Saturday\"],\n  [\"471\u00d76=2826\", \"461\u00d72=922\"],\n  [\"570\u00d79=5130\", \"643\u00d76=3858\"],\n  [\"143\u00d75=715\", \"789\u00d73=2367\"],\n  [\"425\u00d77=2975\", \"115\u00d79=1035\"],\n  [\"657\u00d72=1314\", \"433\u00d73=1299\"],\n  [\"152\u00d78=1216\", \"308\u00d78=2464\"],\n  [\"556\u00d75=2780\", \"494\u00d73=1482\"],\n  [\"881\u00d78=7048\", \"470\u00d73=1410\"],\n  [\"372\u00d77=2604\", \"697\u00d75=3485\"],\n  [\"359\u00d76=2154\", \"177\u00d74=708\"],\n  [\"105\u00d75=525\", \"583\u00d76=3498\"],\n  [\"372\u00d75=1860\", \"598\u00d78=4784\"],\n  [\"754\u00d77=5278\", \"226\u00d79=2034\"],\n  [\"686\u00d79=6174\", \"359\u00d75=1795\"],\n  [\"779\u00d72=1558\", \"651\u00d76=3906\"],\n  [\"563\u00d76=3378\", \"657\u00d76=3942\"],\n  [\"221\u00d75=1105\", \"153\u00d73=459\"],\n  [\"788\u00d72=1576\", \"957\u00d72=1914\"],\n  [\"787\u00d78=6296\", \"614\u00d79=5526\"],\n  [\"940\u00d78=7520\", \"714\u00d77=4998\"],\n  [\"452\u00d78=3616\", \"724\u00d75=3620\"],\n  [\"941\u00d78=7528\", \"951\u00d76=5706\"],\n  [\"908\u00d78=7264\", \"589\u00d73=1767\"],\n  [\"410\u00d74=1640\", \"116\u00d79=1044\"],\n  [\"175\u00d74=700\", \"260\u00d76=1560\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"A\u00d7B=C\" answer cell to the new values\n# from the commit. Each old string is unique in the document, so a\n# Find/Replace over the whole document content is sufficient and leaves\n# the existing run formatting (fonts/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-03-21 Friday\", \"2025-03-22 Saturday\"),\n  @(\"471\u00d76=2826\", \"461\u00d72=922\"),\n  @(\"570\u00d79=5130\", \"643\u00d76=3858\"),\n  @(\"143\u00d75=715\", \"789\u00d73=2367\"),\n  @(\"425\u00d77=2975\", \"115\u00d79=1035\"),\n  @(\"657\u00d72=1314\", \"433\u00d73=1299\"),\n  @(\"152\u00d78=1216\", \"308\u00d78=2464\"),\n  @(\"556\u00d75=2780\", \"494\u00d73=1482\"),\n  @(\"881\u00d78=7048\", \"470\u00d73=1410\"),\n  @(\"372\u00d77=2604\", \"697\u00d75=3485\"),\n  @(\"359\u00d76=2154\", \"177\u00d74=708\"),\n  @(\"105\u00d75=525\", \"583\u00d76=3498\"),\n  @(\"372\u00d75=1860\", \"598\u00d78=4784\"),\n  @(\"754\u00d77=5278\", \"226\u00d79=2034\"),\n  @(\"686\u00d79=6174\", \"359\u00d75=1795\"),\n  @(\"779\u00d72=1558\", \"651\u00d76=3906\"),\n  @(\"563\u00d76=3378\", \"657\u00d76=3942\"),\n  @(\"221\u00d75=1105\", \"153\u00d73=459\"),\n  @(\"788\u00d72=1576\", \"957\u00d72=1914\"),\n  @(\"787\u00d78=6296\", \"614\u00d79=5526\"),\n  @(\"940\u00d78=7520\", \"714\u00d77=4998\"),\n  @(\"452\u00d78=3616\", \"724\u00d75=3620\"),\n  @(\"941\u00d78=7528\", \"951\u00d76=5706\"),\n  @(\"908\u00d78=7264\", \"589\u00d73=1767\"),\n  @(\"410\u00d74=1640\", \"116\u00d79=1044\"),\n  @(\"175\u00d74=700\", \"260\u00d76=1560\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
